# =============================================================================
# 603486-科沃斯.xlsx : add "2022-Q3" fund-holders data
#
#  1. Insert a brand-new worksheet named "2022-Q3" right before "2022-Q2"
#     (pushes every other quarterly tab one slot to the right).
#  2. Populate that new sheet with the fund-holdings table for the quarter.
#  3. Insert a new row at the top of the "总计" (summary) sheet's data and
#     fill it with the 2022-Q3 roll-up figures (持有数量 = 20, 持有市值 = 3.56).
# =============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) New "2022-Q3" worksheet, positioned before "2022-Q2".
# ---------------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($refSheet)
$newSheet.Name = "2022-Q3"

$ws = $newSheet

# Header row B1:H1 - reuse the bold/centered/bordered style already used by
# the other quarter sheets' header row so the cell style index matches.
$headerSrc = $refSheet.Range("B1:H1")
$headerSrc.Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

$ws.Range("A2").Value = 0
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "240008"
$ws.Range("B2").ClearFormats()
$ws.Range("C2").Value = "华宝收益增长混合A"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "8.39"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "93.74"
$ws.Range("E2").ClearFormats()
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "7.87"
$ws.Range("F2").ClearFormats()
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "0.6603"
$ws.Range("G2").ClearFormats()
$ws.Range("H2").Value = 5
$ws.Range("A3").Value = 1
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "159996"
$ws.Range("B3").ClearFormats()
$ws.Range("C3").Value = "国泰中证全指家用电器ETF"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.72"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "97.52"
$ws.Range("E3").ClearFormats()
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "2.79"
$ws.Range("F3").ClearFormats()
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "0.6060"
$ws.Range("G3").ClearFormats()
$ws.Range("H3").Value = 7
$ws.Range("A4").Value = 2
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "010736"
$ws.Range("B4").ClearFormats()
$ws.Range("C4").Value = "易方达沪深300指数精选增强A"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "14.13"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "85.18"
$ws.Range("E4").ClearFormats()
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "3.78"
$ws.Range("F4").ClearFormats()
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "0.5341"
$ws.Range("G4").ClearFormats()
$ws.Range("H4").Value = 6
$ws.Range("A5").Value = 3
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "010020"
$ws.Range("B5").ClearFormats()
$ws.Range("C5").Value = "华夏线上经济主题精选混合"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "15.58"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "85.53"
$ws.Range("E5").ClearFormats()
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "3.39"
$ws.Range("F5").ClearFormats()
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "0.5282"
$ws.Range("G5").ClearFormats()
$ws.Range("H5").Value = 9
$ws.Range("A6").Value = 4
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "011282"
$ws.Range("B6").ClearFormats()
$ws.Range("C6").Value = "华夏消费龙头混合A"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "14.12"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "86.27"
$ws.Range("E6").ClearFormats()
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "3.46"
$ws.Range("F6").ClearFormats()
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "0.4886"
$ws.Range("G6").ClearFormats()
$ws.Range("H6").Value = 10
$ws.Range("A7").Value = 5
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "560880"
$ws.Range("B7").ClearFormats()
$ws.Range("C7").Value = "广发中证全指家用电器ETF"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "9.25"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "99.30"
$ws.Range("E7").ClearFormats()
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "2.91"
$ws.Range("F7").ClearFormats()
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "0.2692"
$ws.Range("G7").ClearFormats()
$ws.Range("H7").Value = 6
$ws.Range("A8").Value = 6
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "010737"
$ws.Range("B8").ClearFormats()
$ws.Range("C8").Value = "易方达沪深300指数精选增强C"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.64"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "85.18"
$ws.Range("E8").ClearFormats()
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "3.78"
$ws.Range("F8").ClearFormats()
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "0.1376"
$ws.Range("G8").ClearFormats()
$ws.Range("H8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "010692"
$ws.Range("B9").ClearFormats()
$ws.Range("C9").Value = "华夏核心价值混合A"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.83"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "89.93"
$ws.Range("E9").ClearFormats()
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "3.50"
$ws.Range("F9").ClearFormats()
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "0.0640"
$ws.Range("G9").ClearFormats()
$ws.Range("H9").Value = 10
$ws.Range("A10").Value = 8
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "011283"
$ws.Range("B10").ClearFormats()
$ws.Range("C10").Value = "华夏消费龙头混合C"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.25"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "86.27"
$ws.Range("E10").ClearFormats()
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "3.46"
$ws.Range("F10").ClearFormats()
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "0.0432"
$ws.Range("G10").ClearFormats()
$ws.Range("H10").Value = 10
$ws.Range("A11").Value = 9
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "562500"
$ws.Range("B11").ClearFormats()
$ws.Range("C11").Value = "华夏中证机器人ETF"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.57"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "99.51"
$ws.Range("E11").ClearFormats()
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "2.51"
$ws.Range("F11").ClearFormats()
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "0.0394"
$ws.Range("G11").ClearFormats()
$ws.Range("H11").Value = 8
$ws.Range("A12").Value = 10
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "561120"
$ws.Range("B12").ClearFormats()
$ws.Range("C12").Value = "富国中证全指家用电器ETF"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.27"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "99.27"
$ws.Range("E12").ClearFormats()
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "2.91"
$ws.Range("F12").ClearFormats()
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "0.0370"
$ws.Range("G12").ClearFormats()
$ws.Range("H12").Value = 6
$ws.Range("A13").Value = 11
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "159770"
$ws.Range("B13").ClearFormats()
$ws.Range("C13").Value = "天弘中证机器人ETF"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.04"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "99.68"
$ws.Range("E13").ClearFormats()
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "2.50"
$ws.Range("F13").ClearFormats()
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "0.0260"
$ws.Range("G13").ClearFormats()
$ws.Range("H13").Value = 8
$ws.Range("A14").Value = 12
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "012461"
$ws.Range("B14").ClearFormats()
$ws.Range("C14").Value = "西藏东财国证龙头家电指数A"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.66"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "94.59"
$ws.Range("E14").ClearFormats()
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "3.81"
$ws.Range("F14").ClearFormats()
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0.0251"
$ws.Range("G14").ClearFormats()
$ws.Range("H14").Value = 7
$ws.Range("A15").Value = 13
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "012462"
$ws.Range("B15").ClearFormats()
$ws.Range("C15").Value = "西藏东财国证龙头家电指数C"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.59"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "94.59"
$ws.Range("E15").ClearFormats()
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "3.81"
$ws.Range("F15").ClearFormats()
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "0.0225"
$ws.Range("G15").ClearFormats()
$ws.Range("H15").Value = 7
$ws.Range("A16").Value = 14
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "562360"
$ws.Range("B16").ClearFormats()
$ws.Range("C16").Value = "银华中证机器人ETF"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.77"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "97.23"
$ws.Range("E16").ClearFormats()
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "2.45"
$ws.Range("F16").ClearFormats()
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "0.0189"
$ws.Range("G16").ClearFormats()
$ws.Range("H16").Value = 8
$ws.Range("A17").Value = 15
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "010693"
$ws.Range("B17").ClearFormats()
$ws.Range("C17").Value = "华夏核心价值混合C"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.54"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "89.93"
$ws.Range("E17").ClearFormats()
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = "3.50"
$ws.Range("F17").ClearFormats()
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "0.0189"
$ws.Range("G17").ClearFormats()
$ws.Range("H17").Value = 10
$ws.Range("A18").Value = 16
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "013054"
$ws.Range("B18").ClearFormats()
$ws.Range("C18").Value = "天弘国证龙头家电指数C"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.49"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "94.83"
$ws.Range("E18").ClearFormats()
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = "3.82"
$ws.Range("F18").ClearFormats()
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "0.0187"
$ws.Range("G18").ClearFormats()
$ws.Range("H18").Value = 7
$ws.Range("A19").Value = 17
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "159730"
$ws.Range("B19").ClearFormats()
$ws.Range("C19").Value = "博时国证龙头家电ETF"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.45"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "98.01"
$ws.Range("E19").ClearFormats()
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = "3.96"
$ws.Range("F19").ClearFormats()
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "0.0178"
$ws.Range("G19").ClearFormats()
$ws.Range("H19").Value = 7
$ws.Range("A20").Value = 18
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "013053"
$ws.Range("B20").ClearFormats()
$ws.Range("C20").Value = "天弘国证龙头家电指数A"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.16"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "94.83"
$ws.Range("E20").ClearFormats()
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = "3.82"
$ws.Range("F20").ClearFormats()
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "0.0061"
$ws.Range("G20").ClearFormats()
$ws.Range("H20").Value = 7
$ws.Range("A21").Value = 19
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "015573"
$ws.Range("B21").ClearFormats()
$ws.Range("C21").Value = "华宝收益增长混合C"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.04"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "93.74"
$ws.Range("E21").ClearFormats()
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = "7.87"
$ws.Range("F21").ClearFormats()
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "0.0031"
$ws.Range("G21").ClearFormats()
$ws.Range("H21").Value = 5
# Column-A row-index cells (A2:A21) use the same bold/centered/bordered style.
$colASrc = $refSheet.Range("A2")
$colASrc.Copy()
$ws.Range("A2:A21").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) "总计" (summary) sheet: insert the 2022-Q3 roll-up as the new row 2.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 20
$summary.Range("D2").Value = 3.56
$summary.Range("B2:D2").ClearFormats()

$sumColASrc = $summary.Range("A3")
$sumColASrc.Copy()
$summary.Range("A2").PasteSpecial(-4122)
